# materialprofiles.xlsx — PLA purge-bubble update + new HIPS profile row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 (rigid.ink HIPS): print settings profile + notes filled in ---
# Print settings column (B10) changes from "0.15 QUALITY MK3" to "rigid.ink HIPS"
# and loses its left-border formatting (becomes unstyled like B9/B11).
$ws.Cells.Item(10, 2).Value = "rigid.ink HIPS"
$ws.Cells.Item(10, 2).Borders.Item(7).LineStyle = -4142

# Notes column (E10) gets a new comment about print quality.
$ws.Cells.Item(10, 5).Value = "Not quite happy with print quality, oozey and stringy"

# --- Row 8 (rigid.ink PLA): printer profile switched to the "purgebubble" variant ---
$ws.Cells.Item(8, 4).Value = "Original Prusa i3 MK3 purgebubble"

# --- Column D widened to fit the new, longer printer-profile text ---
$ws.Columns.Item(4).ColumnWidth = 31.25

# --- Selection left where the author's cursor ended up ---
$ws.Range("E15").Select()
